# Remove footnote-reference markers like " [1]" / " [5, 6]" from vaccine
# names and descriptions, and collapse cells whose text was wrapped onto a
# second line (an embedded newline) back into a single line (joined with a
# single space), across every worksheet in the workbook.
#
# e.g. "DTaP [1]"                              -> "DTaP "
#      "Hepatitis B [5]\nPediatric/Adolescent"  -> "Hepatitis B  Pediatric/Adolescent"
#      "Recombivax\nHB"                         -> "Recombivax HB"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowOffset = $used.Row
    $colOffset = $used.Column
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($rowOffset + $r, $colOffset + $c)
            $val = $cell.Value2

            if ($val -is [string]) {
                $newVal = $val -replace '\[[0-9]+(?:,\s*[0-9]+)*\]', ''
                $newVal = $newVal -replace "`r`n", ' '
                $newVal = $newVal -replace "`n", ' '

                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
